# recreation of merges, pivot tables, charts, and conditional formatting of
# various dataframes -- this particular edit fixes a pandas merge-suffix
# mixup: the seas_id_x/player_id/seas_id_y/season_ending_year_y headers
# (and their K/L/N/O column data) were produced by merging on the wrong
# keys, so columns K, L, N, O need their header labels corrected and the
# N/O column values (which had been swapped) put back in the right place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row fixes (row 1) ---------------------------------------------
$ws.Range("K1").Value = "seas_id"
$ws.Range("L1").Value = "player_id_x"
$ws.Range("N1").Value = "season_ending_year_y"
$ws.Range("O1").Value = "player_id_y"

# --- Data fixes: columns N and O were swapped for rows 2-11 ---------------
# N should hold the (text) "season_ending_year_y" value that used to sit in
# O, and O should hold the (numeric) "player_id_y" value that used to sit
# in N. We stage the text values through a scratch cell formatted as text
# so they land as genuine shared-string cells (matching how the season
# year values were already stored), then paste values-only into N so no
# stray number formatting sticks to the destination cells.

$scratch = $ws.Range("ZZ1")

function Set-TextValue($cell, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.NumberFormat = "General"
    $scratch.Clear()
}

Set-TextValue $ws.Range("N2") "1976"
$ws.Range("O2").Value = 2932

Set-TextValue $ws.Range("N3") "1975"
$ws.Range("O3").Value = 2932

Set-TextValue $ws.Range("N4") "1975"
$ws.Range("O4").Value = 1906

Set-TextValue $ws.Range("N5") "1974"
$ws.Range("O5").Value = 2932

Set-TextValue $ws.Range("N6") "1973"
$ws.Range("O6").Value = 426

Set-TextValue $ws.Range("N7") "1972"
$ws.Range("O7").Value = 293

Set-TextValue $ws.Range("N8") "1971"
$ws.Range("O8").Value = 3599

Set-TextValue $ws.Range("N9") "1970"
$ws.Range("O9").Value = 4639

Set-TextValue $ws.Range("N10") "1969"
$ws.Range("O10").Value = 3599

Set-TextValue $ws.Range("N11") "1968"
$ws.Range("O11").Value = 1008
